# Apply updated figures to the "Crédito disponível - Centraliza" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N11").Value = 488476.82
$ws.Range("O11").Value = 462793.92

$ws.Range("O13").Value = 13488.83

$ws.Range("K17").Value = 38733.48

$ws.Range("K22").Value = 1800

$ws.Range("K23").Value = 22884.17

$ws.Range("K25").Value = 26301

$ws.Range("M26").Value = 126410
